$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header row + data rows).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Rename header cells: "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410"
# (the "diff" header, and any other header without that suffix, is left untouched).
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        $newVal = $v -replace "_old$", "_FV2404"
        $newVal = $newVal -replace "_new$", "_FV2410"
        if ($newVal -ne $v) {
            $cell.Value = $newVal
        }
    }
}

# Turn the data range into an Excel Table ("Table1") with a header row and autofilter.
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $headerRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, top-left of the scrollable area is A2).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
